# Auto-generated update of NATMI LR-pair TPM-derived metrics (Col1a2-Itgb1)
# Commit: "update scripts wuth new tpm"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 17.73076433333334
$ws.Range("H2").Value = 53.19229300000001
$ws.Range("I2").Value = 0.004631884691211661
$ws.Range("J2").Value = 0.00463188469121166
$ws.Range("M2").Value = 168.1098273333333
$ws.Range("N2").Value = 504.329482
$ws.Range("O2").Value = 0.2984182258032519
$ws.Range("P2").Value = 0.298418225803252
$ws.Range("Q2").Value = 2980.715730564692
$ws.Range("R2").Value = 26826.44157508223
$ws.Range("S2").Value = 0.001382238811676627
$ws.Range("T2").Value = 0.001382238811676627
$ws.Range("G3").Value = 17.73076433333334
$ws.Range("H3").Value = 53.19229300000001
$ws.Range("I3").Value = 0.004631884691211661
$ws.Range("J3").Value = 0.00463188469121166
$ws.Range("O3").Value = 0.2893586437755394
$ws.Range("P3").Value = 0.2893586437755394
$ws.Range("Q3").Value = 2890.225149469462
$ws.Range("R3").Value = 26012.02634522515
$ws.Range("S3").Value = 0.001340275872373689
$ws.Range("T3").Value = 0.001340275872373689
$ws.Range("G4").Value = 17.73076433333334
$ws.Range("H4").Value = 53.19229300000001
$ws.Range("I4").Value = 0.004631884691211661
$ws.Range("J4").Value = 0.00463188469121166
$ws.Range("M4").Value = 165.99353
$ws.Range("N4").Value = 497.98059
$ws.Range("O4").Value = 0.294661504941043
$ws.Range("P4").Value = 0.294661504941043
$ws.Range("Q4").Value = 2943.192161288097
$ws.Range("R4").Value = 26488.72945159287
$ws.Range("S4").Value = 0.001364838113825806
$ws.Range("T4").Value = 0.001364838113825806
$ws.Range("G5").Value = 17.73076433333334
$ws.Range("H5").Value = 53.19229300000001
$ws.Range("I5").Value = 0.004631884691211661
$ws.Range("J5").Value = 0.00463188469121166
$ws.Range("M5").Value = 66.22673433333334
$ws.Range("N5").Value = 198.680203
$ws.Range("O5").Value = 0.1175616254801657
$ws.Range("P5").Value = 0.1175616254801657
$ws.Range("Q5").Value = 1174.250619030609
$ws.Range("R5").Value = 10568.25557127548
$ws.Range("S5").Value = 0.000544531893335538
$ws.Range("T5").Value = 0.0005445318933355379
$ws.Range("I6").Value = 0.9353873458333681
$ws.Range("J6").Value = 0.935387345833368
$ws.Range("M6").Value = 168.1098273333333
$ws.Range("N6").Value = 504.329482
$ws.Range("O6").Value = 0.2984182258032519
$ws.Range("P6").Value = 0.298418225803252
$ws.Range("Q6").Value = 601941.5339044909
$ws.Range("R6").Value = 5417473.805140418
$ws.Range("S6").Value = 0.2791366321824066
$ws.Range("T6").Value = 0.2791366321824066
$ws.Range("I7").Value = 0.9353873458333681
$ws.Range("J7").Value = 0.935387345833368
$ws.Range("O7").Value = 0.2893586437755394
$ws.Range("P7").Value = 0.2893586437755394
$ws.Range("R7").Value = 5253006.476817743
$ws.Range("S7").Value = 0.2706624137951449
$ws.Range("T7").Value = 0.2706624137951448
$ws.Range("I8").Value = 0.9353873458333681
$ws.Range("J8").Value = 0.935387345833368
$ws.Range("M8").Value = 165.99353
$ws.Range("N8").Value = 497.98059
$ws.Range("O8").Value = 0.294661504941043
$ws.Range("P8").Value = 0.294661504941043
$ws.Range("Q8").Value = 594363.8254312156
$ws.Range("R8").Value = 5349274.42888094
$ws.Range("S8").Value = 0.2756226430260681
$ws.Range("T8").Value = 0.275622643026068
$ws.Range("I9").Value = 0.9353873458333681
$ws.Range("J9").Value = 0.935387345833368
$ws.Range("M9").Value = 66.22673433333334
$ws.Range("N9").Value = 198.680203
$ws.Range("O9").Value = 0.1175616254801657
$ws.Range("P9").Value = 0.1175616254801657
$ws.Range("Q9").Value = 237134.3941187155
$ws.Range("R9").Value = 2134209.54706844
$ws.Range("S9").Value = 0.1099656568297486
$ws.Range("T9").Value = 0.1099656568297486
$ws.Range("G10").Value = 227.2177583333333
$ws.Range("H10").Value = 681.653275
$ws.Range("I10").Value = 0.0593570833501536
$ws.Range("J10").Value = 0.05935708335015359
$ws.Range("M10").Value = 168.1098273333333
$ws.Range("N10").Value = 504.329482
$ws.Range("O10").Value = 0.2984182258032519
$ws.Range("P10").Value = 0.298418225803252
$ws.Range("Q10").Value = 38197.53812048372
$ws.Range("R10").Value = 343777.8430843535
$ws.Range("S10").Value = 0.01771323550220858
$ws.Range("T10").Value = 0.01771323550220858
$ws.Range("G11").Value = 227.2177583333333
$ws.Range("H11").Value = 681.653275
$ws.Range("I11").Value = 0.0593570833501536
$ws.Range("J11").Value = 0.05935708335015359
$ws.Range("O11").Value = 0.2893586437755394
$ws.Range("P11").Value = 0.2893586437755394
$ws.Range("Q11").Value = 37037.91146253505
$ws.Range("R11").Value = 333341.2031628154
$ws.Range("S11").Value = 0.0171754851366721
$ws.Range("T11").Value = 0.01717548513667209
$ws.Range("G12").Value = 227.2177583333333
$ws.Range("H12").Value = 681.653275
$ws.Range("I12").Value = 0.0593570833501536
$ws.Range("J12").Value = 0.05935708335015359
$ws.Range("M12").Value = 165.99353
$ws.Range("N12").Value = 497.98059
$ws.Range("O12").Value = 0.294661504941043
$ws.Range("P12").Value = 0.294661504941043
$ws.Range("Q12").Value = 37716.67778443691
$ws.Range("R12").Value = 339450.1000599323
$ws.Range("S12").Value = 0.01749024750886718
$ws.Range("T12").Value = 0.01749024750886718
$ws.Range("G13").Value = 227.2177583333333
$ws.Range("H13").Value = 681.653275
$ws.Range("I13").Value = 0.0593570833501536
$ws.Range("J13").Value = 0.05935708335015359
$ws.Range("M13").Value = 66.22673433333334
$ws.Range("N13").Value = 198.680203
$ws.Range("O13").Value = 0.1175616254801657
$ws.Range("P13").Value = 0.1175616254801657
$ws.Range("Q13").Value = 15047.8901169572
$ws.Range("R13").Value = 135431.0110526148
$ws.Range("S13").Value = 0.006978115202405735
$ws.Range("T13").Value = 0.006978115202405734
$ws.Range("G14").Value = 2.387458333333333
$ws.Range("H14").Value = 7.162374999999999
$ws.Range("I14").Value = 0.0006236861252666267
$ws.Range("J14").Value = 0.0006236861252666266
$ws.Range("M14").Value = 168.1098273333333
$ws.Range("N14").Value = 504.329482
$ws.Range("O14").Value = 0.2984182258032519
$ws.Range("P14").Value = 0.298418225803252
$ws.Range("Q14").Value = 401.3552081821944
$ws.Range("R14").Value = 3612.196873639749
$ws.Range("S14").Value = 0.0001861193069601715
$ws.Range("T14").Value = 0.0001861193069601715
$ws.Range("G15").Value = 2.387458333333333
$ws.Range("H15").Value = 7.162374999999999
$ws.Range("I15").Value = 0.0006236861252666267
$ws.Range("J15").Value = 0.0006236861252666266
$ws.Range("O15").Value = 0.2893586437755394
$ws.Range("P15").Value = 0.2893586437755394
$ws.Range("Q15").Value = 389.1705957276805
$ws.Range("R15").Value = 3502.535361549124
$ws.Range("S15").Value = 0.0001804689713487723
$ws.Range("T15").Value = 0.0001804689713487723
$ws.Range("G16").Value = 2.387458333333333
$ws.Range("H16").Value = 7.162374999999999
$ws.Range("I16").Value = 0.0006236861252666267
$ws.Range("J16").Value = 0.0006236861252666266
$ws.Range("M16").Value = 165.99353
$ws.Range("N16").Value = 497.98059
$ws.Range("O16").Value = 0.294661504941043
$ws.Range("P16").Value = 0.294661504941043
$ws.Range("Q16").Value = 396.3026364779166
$ws.Range("R16").Value = 3566.72372830125
$ws.Range("S16").Value = 0.0001837762922819121
$ws.Range("T16").Value = 0.000183776292281912
$ws.Range("G17").Value = 2.387458333333333
$ws.Range("H17").Value = 7.162374999999999
$ws.Range("I17").Value = 0.0006236861252666267
$ws.Range("J17").Value = 0.0006236861252666266
$ws.Range("M17").Value = 66.22673433333334
$ws.Range("N17").Value = 198.680203
$ws.Range("O17").Value = 0.1175616254801657
$ws.Range("P17").Value = 0.1175616254801657
$ws.Range("Q17").Value = 158.1135687735694
$ws.Range("R17").Value = 1423.022118962125
$ws.Range("S17").Value = 0.00007332155467577086
$ws.Range("T17").Value = 0.00007332155467577085
